$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values updated
$ws.Range("B3").Value = 0.9976881759438693
$ws.Range("C3").Value = 0.9977081588213691
$ws.Range("D3").Value = 0.9975001834993759

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9973570029753196
$ws.Range("C4").Value = 0.9974613023883855
$ws.Range("D4").Value = 0.9974807931223996

# Row 5: AdaBoostRegressor -> MLPRegressor, values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9979699599778442
$ws.Range("C5").Value = 0.9980309405920186
$ws.Range("D5").Value = 0.9979966486923662
